$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "562 WILLIAM BERCZY BLVD MARKHAM ON L6C2P7",
    "MARKHAM WILLIAM BERCZY BLVD ON L6C2P7 562",
    "WILLIAM BERCZY BLVD 562 ON L6C2P7 MARKHAM",
    "WILLIAM BERCZY BLVD MARKHAM ON L6C2P7 562",
    "562 MARKHAM ON L6C2P7 WILLIAM BERCZY BLVD",
    "MARKHAM 562 ON L6C2P7 WILLIAM BERCZY BLVD"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("J5").Select() | Out-Null
